# Generate Report for Handback
# Update "Latest Handback DateTime" (column K) for the f94e23ab file (row 2)
# on both the zh-cn and de-de localization status sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("K2").Value = "2016-10-27 10:12:23"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-10-27 10:12:39"
